{"js": "// Remove the standalone \"Meta description: ...\" paragraph from the top of\n// the document, and add a \"Play Diamond Mystic Slot Game for Free | Slot\n// Review\" (bold) paragraph plus the meta-description text (italic) right\n// before the final \"Create an image...\" paragraph, replacing that\n// paragraph's old image-prompt text with the meta description text.\n\nconst body = context.document.body;\n\nconst metaDescriptionText =\n  \"Read our Diamond Mystic slot game review to play for free. Get info on impressive payouts, simplicity, low betting range, and limited features.\";\nconst titleText = \"Play Diamond Mystic Slot Game for Free | Slot Review\";\n\n// --- Step 1: delete the \"Meta description\" paragraph right after the title ---\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Meta description\") === 0) {\n    paras.items[i].delete();\n    break;\n  }\n}\nawait context.sync();\n\n// --- Step 2: locate the final paragraph (the old image-prompt paragraph) ---\nconst paras2 = body.paragraphs;\nparas2.load(\"text\");\nawait context.sync();\n\nconst finalIndex = paras2.items.length - 1;\nconst finalPara = paras2.items[finalIndex];\n\n// Insert a new, empty paragraph right before it.\nconst insertionRange = finalPara.getRange(\"Start\");\ninsertionRange.insertParagraph(\"\", \"Before\");\nawait context.sync();\n\n// Re-fetch the paragraph collection so the freshly inserted paragraph (and\n// the final paragraph after it) are addressed cleanly.\nconst paras3 = body.paragraphs;\nparas3.load(\"text\");\nawait context.sync();\n\nconst newParaIndex = paras3.items.length - 2;\nconst finalParaIndex = paras3.items.length - 1;\nconst newPara = paras3.items[newParaIndex];\nconst finalPara2 = paras3.items[finalParaIndex];\n\n// Fill the new paragraph with exact OOXML so it gets clean bold formatting\n// with no italic/style inherited from its neighboring paragraphs.\nconst newParaRange = newPara.getRange(\"Whole\");\nconst titleOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>\" +\n  titleText +\n  \"</w:t></w:r></w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\nnewParaRange.insertOoxml(titleOoxml, \"Replace\");\nawait context.sync();\n\n// Replace the old image-prompt text in the final paragraph with the meta\n// description text, keeping its existing italic run formatting.\nconst finalRange = finalPara2.getRange(\"Whole\");\nfinalRange.insertText(metaDescriptionText, \"Replace\");\nawait context.sync();\n", "ps1": "# Remove the standalone \"Meta description: ...\" paragraph from the top of\n# the document, and add a \"Play Diamond Mystic Slot Game for Free | Slot\n# Review\" (bold) paragraph plus the meta-description text (italic) right\n# before the final \"Create an image...\" paragraph, replacing that\n# paragraph's old image-prompt text with the meta description text.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Diamond Mystic Slot Game for Free | Slot Review\"\n$metaDescriptionText = \"Read our Diamond Mystic slot game review to play for free. Get info on impressive payouts, simplicity, low betting range, and limited features.\"\n\n# --- Step 1: delete the \"Meta description\" paragraph right after the title ---\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"Meta description\")) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- Step 2: insert a new empty paragraph right before the final paragraph\n#     (the old image-prompt paragraph) ---\n$count2 = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count2)\n$lastPara.Range.InsertParagraphBefore()\n\n# Fill the newly inserted paragraph with exact OOXML so it gets clean bold\n# formatting with no italic/style inherited from its neighboring paragraphs.\n$count3 = $d.Paragraphs.Count\n$newPara = $d.Paragraphs.Item($count3 - 1)\n$newRange = $newPara.Range\n$titleOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$newRange.InsertXML($titleOoxml)\n\n# --- Step 3: replace the old image-prompt text in the final paragraph with\n#     the meta description text, keeping its existing italic run formatting ---\n$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$finalRange = $finalPara.Range\n$finalRange.MoveEnd(1, -1)\n$finalRange.Text = $metaDescriptionText\n"}
